# Modules.xlsx rework:
#   - the "CodeModule" column (was column F) becomes column A; every other
#     column (ModuleName, ElementName1, ElementName2, Dept_Attachement,
#     Coordinator) shifts one slot to the right (A->B, B->C, ... E->F).
#   - the CodeModule values are renamed from AP2x to AP1x.
#   - the active selection moves to A3 (and any leftover horizontal scroll
#     position is cleared).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 10
$lastCol = 6

# Snapshot the current grid (by display text) before we start overwriting
# cells, since writes would otherwise clobber values we still need to read.
$data = @()
for ($r = 1; $r -le $lastRow; $r++) {
    $row = @()
    for ($c = 1; $c -le $lastCol; $c++) {
        $row += $ws.Cells.Item($r, $c).Text
    }
    $data += ,$row
}

# Write the data back shifted: new column 1 <- old column 6 (CodeModule),
# new column 2 <- old column 1, new column 3 <- old column 2, etc.
for ($r = 1; $r -le $lastRow; $r++) {
    $old = $data[$r - 1]
    $ws.Cells.Item($r, 1).Value = $old[5]
    $ws.Cells.Item($r, 2).Value = $old[0]
    $ws.Cells.Item($r, 3).Value = $old[1]
    $ws.Cells.Item($r, 4).Value = $old[2]
    $ws.Cells.Item($r, 5).Value = $old[3]
    $ws.Cells.Item($r, 6).Value = $old[4]
}

# Rename the CodeModule values (now in column A) from AP2x to AP1x.
for ($r = 2; $r -le $lastRow; $r++) {
    $code = $ws.Cells.Item($r, 1).Text
    if ($code -like "AP2*") {
        $ws.Cells.Item($r, 1).Value = "AP1" + $code.Substring(3)
    }
}

# Move the selection/active cell to A3 and reset any scrolled view.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("A3").Select() | Out-Null
